$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.457.77"
$ws.Range("E2").Value = "  -1.47%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.357.32"
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.32"
$ws.Range("E5").Value = "  +5.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.25"
$ws.Range("E6").Value = "  -7.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.636"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.620"
$ws.Range("E9").Value = "  -0.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.13"
$ws.Range("E10").Value = "  -7.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0918"
$ws.Range("E11").Value = "  -2.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.41"
$ws.Range("E12").Value = "  -6.44%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.01"
$ws.Range("E13").Value = "  -5.13%  "
$ws.Range("E14").Value = "  +0.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.22"
$ws.Range("E15").Value = "  +0.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.717.67"
$ws.Range("E16").Value = "  +0.54%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.364.28"
$ws.Range("E17").Value = "  -2.70%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.578.49"
$ws.Range("E18").Value = "  -1.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.58"
$ws.Range("E19").Value = "  +4.55%  "
$ws.Range("E20").Value = "  -1.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "75.11"
$ws.Range("E21").Value = "  -0.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.70"
$ws.Range("E22").Value = "  +7.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "275.62"
$ws.Range("E23").Value = "  +9.42%  "
$ws.Range("E24").Value = "  -8.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.69"
$ws.Range("E25").Value = "  +8.12%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.47"
$ws.Range("E27").Value = "  -4.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.80"
$ws.Range("E28").Value = "  +6.07%  "
$ws.Range("E29").Value = "  -1.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "173.85"
$ws.Range("E30").Value = "  +0.58%  "
$ws.Range("E31").Value = "  -2.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "35.50"
$ws.Range("E32").Value = "  -8.22%  "
$ws.Range("E33").Value = "  -0.82%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.93"
$ws.Range("E34").Value = "  +2.46%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.133"
$ws.Range("E35").Value = "  +1.77%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.61"
$ws.Range("E36").Value = "  -6.91%  "
$ws.Range("E37").Value = "  -5.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.88"
$ws.Range("E38").Value = "  -6.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.88"
$ws.Range("E39").Value = "  +4.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.105"
$ws.Range("E40").Value = "  +0.60%  "
$ws.Range("E41").Value = "  +0.84%  "
$ws.Range("E42").Value = "  -2.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "68.81"
$ws.Range("E43").Value = "  -4.63%  "
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "90.30"
$ws.Range("E45").Value = "  +35.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "115.34"
$ws.Range("E46").Value = "  +4.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.90"
$ws.Range("E47").Value = "  -5.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.46"
$ws.Range("E48").Value = "  -4.06%  "
$ws.Range("E49").Value = "  -1.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.586.46"
$ws.Range("E50").Value = "  +6.39%  "
$ws.Range("E51").Value = "  -2.69%  "
